# Commit: "Removed source image from orthogonal image as already contained
# in image superclass. Removed wrong if_absent value"
#
# 1) OrthogonalImage: drop the now-redundant leading "source_image" column
#    (column A) - it is already provided by the Image superclass. Every
#    other column (source_roi, axis, ...) shifts one position to the left,
#    and the "axis" list-validation (originally on column C) shifts with it
#    to column B.
# 2) UserExperimentKeyMeasurements: insert a new "saturated_channels" column
#    right before "table_data".

$wb = $excel.ActiveWorkbook

# --- OrthogonalImage: remove the redundant "source_image" column (A) ---
$wsOrtho = $wb.Worksheets.Item("OrthogonalImage")
$wsOrtho.Range("A1").EntireColumn.Delete()

# --- UserExperimentKeyMeasurements: insert "saturated_channels" before table_data (D) ---
$wsUEKM = $wb.Worksheets.Item("UserExperimentKeyMeasurements")
$wsUEKM.Range("D1").EntireColumn.Insert()
$wsUEKM.Range("D1").Value = "saturated_channels"
